$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 25050150
$ws.Range("I53").Value = 55555636
$ws.Range("J53").Value = 91115.63
$ws.Range("K53").Value = 55555636
$ws.Range("L53").Value = 91115.63
$ws.Range("M53").Value = -55554999
$ws.Range("N53").Value = -92389.63

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2018.619
$ws.Range("I132").Value = 2199.7058
$ws.Range("J132").Value = 1249
$ws.Range("K132").Value = 6599.117400000001
$ws.Range("L132").Value = 3747
$ws.Range("M132").Value = -4069.117400000001
$ws.Range("N132").Value = -8807

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6605.783
$ws.Range("I32").Value = 4595.255
$ws.Range("K32").Value = 4595.255
$ws.Range("M32").Value = -4308.255

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6302.423
$ws.Range("I61").Value = 6968.15
$ws.Range("J61").Value = 4083.3333
$ws.Range("K61").Value = 6968.15
$ws.Range("L61").Value = 4083.3333
$ws.Range("M61").Value = -6756.15
$ws.Range("N61").Value = -4507.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 250003460
$ws.Range("I63").Value = 250003460
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 250003460
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -250002774
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 250003460
$ws.Range("I66").Value = 250003460
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 1250017300
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -1250013868
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3226.25
$ws.Range("I132").Value = 1780.8334
$ws.Range("K132").Value = 5342.5002
$ws.Range("M132").Value = -2812.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6302.423
$ws.Range("I136").Value = 6968.15
$ws.Range("J136").Value = 4083.3333
$ws.Range("K136").Value = 20904.45
$ws.Range("L136").Value = 12249.9999
$ws.Range("M136").Value = -18354.45
$ws.Range("N136").Value = -17349.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1288.125
$ws.Range("I94").Value = 788.2273
$ws.Range("J94").Value = 2387.9
$ws.Range("K94").Value = 788.2273
$ws.Range("L94").Value = 2387.9
$ws.Range("M94").Value = -337.2273
$ws.Range("N94").Value = -3289.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2273.818
$ws.Range("I16").Value = 2119.8
$ws.Range("J16").Value = 2402.1667
$ws.Range("K16").Value = 2119.8
$ws.Range("L16").Value = 2402.1667
$ws.Range("M16").Value = -1832.8
$ws.Range("N16").Value = -2976.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4113.9653
$ws.Range("I31").Value = 1706.579
$ws.Range("J31").Value = 8688
$ws.Range("K31").Value = 1706.579
$ws.Range("L31").Value = 8688
$ws.Range("M31").Value = -1411.579
$ws.Range("N31").Value = -9278

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4113.9653
$ws.Range("I34").Value = 1706.579
$ws.Range("J34").Value = 8688
$ws.Range("K34").Value = 1706.579
$ws.Range("L34").Value = 8688
$ws.Range("M34").Value = -1504.579
$ws.Range("N34").Value = -9092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4740.7144
$ws.Range("I99").Value = 5624.5454
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 5624.5454
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -4126.5454
$ws.Range("N99").Value = -4496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1665.0555
$ws.Range("I105").Value = 1947.1
$ws.Range("J105").Value = 1312.5
$ws.Range("K105").Value = 1947.1
$ws.Range("L105").Value = 1312.5
$ws.Range("M105").Value = -200.0999999999999
$ws.Range("N105").Value = -4806.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1699.375
$ws.Range("I107").Value = 942.1429000000001
$ws.Range("K107").Value = 942.1429000000001
$ws.Range("M107").Value = 977.8570999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2273.818
$ws.Range("I113").Value = 2119.8
$ws.Range("J113").Value = 2402.1667
$ws.Range("K113").Value = 2119.8
$ws.Range("L113").Value = 2402.1667
$ws.Range("M113").Value = 50.19999999999982
$ws.Range("N113").Value = -6742.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4740.7144
$ws.Range("I126").Value = 5624.5454
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 16873.6362
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -14403.6362
$ws.Range("N126").Value = -9440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3468.9
$ws.Range("I134").Value = 3982.3333
$ws.Range("J134").Value = 2270.889
$ws.Range("K134").Value = 11946.9999
$ws.Range("L134").Value = 6812.667
$ws.Range("M134").Value = -9411.999899999999
$ws.Range("N134").Value = -11882.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 11111355
$ws.Range("I98").Value = 220.6
$ws.Range("J98").Value = 25000272
$ws.Range("K98").Value = 661.8
$ws.Range("L98").Value = 75000816
$ws.Range("M98").Value = 836.2
$ws.Range("N98").Value = -75003812

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 833868.0600000001
$ws.Range("J113").Value = 1250507.9
$ws.Range("L113").Value = 3751523.7
$ws.Range("N113").Value = -3755863.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2318.818
$ws.Range("J132").Value = 3750.6667
$ws.Range("L132").Value = 33756.0003
$ws.Range("N132").Value = -38816.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5333.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6020.52
$ws.Range("I70").Value = 6265
$ws.Range("J70").Value = 5501
$ws.Range("K70").Value = 6265
$ws.Range("L70").Value = 5501
$ws.Range("M70").Value = -5995
$ws.Range("N70").Value = -6041

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6020.52
$ws.Range("I73").Value = 6265
$ws.Range("J73").Value = 5501
$ws.Range("K73").Value = 6265
$ws.Range("L73").Value = 5501
$ws.Range("M73").Value = -5329
$ws.Range("N73").Value = -7373

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 459152.6
$ws.Range("I102").Value = 848376.75
$ws.Range("J102").Value = 1241.8235
$ws.Range("K102").Value = 848376.75
$ws.Range("L102").Value = 1241.8235
$ws.Range("M102").Value = -846754.75
$ws.Range("N102").Value = -4485.8235

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 33156476
$ws.Range("I122").Value = 64237836
$ws.Range("J122").Value = 3025.8667
$ws.Range("K122").Value = 192713508
$ws.Range("L122").Value = 9077.6001
$ws.Range("M122").Value = -192711058
$ws.Range("N122").Value = -13977.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 45032.824
$ws.Range("I7").Value = 49035.953
$ws.Range("K7").Value = 49035.953
$ws.Range("M7").Value = -48923.953

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 8460.75
$ws.Range("I26").Value = 3333
$ws.Range("J26").Value = 10170
$ws.Range("K26").Value = 3333
$ws.Range("L26").Value = 10170
$ws.Range("M26").Value = -3038
$ws.Range("N26").Value = -10760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 9006.223
$ws.Range("I29").Value = 5016
$ws.Range("K29").Value = 5016
$ws.Range("M29").Value = -4721

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 9393.909
$ws.Range("I33").Value = 3333
$ws.Range("J33").Value = 10000
$ws.Range("K33").Value = 3333
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = -3043
$ws.Range("N33").Value = -10580

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 62502916
$ws.Range("I40").Value = 90910790
$ws.Range("J40").Value = 5589
$ws.Range("K40").Value = 90910790
$ws.Range("L40").Value = 5589
$ws.Range("M40").Value = -90910654
$ws.Range("N40").Value = -5861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 959.6429000000001
$ws.Range("I93").Value = 703.8889
$ws.Range("K93").Value = 703.8889
$ws.Range("M93").Value = 544.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 45032.824
$ws.Range("I126").Value = 49035.953
$ws.Range("K126").Value = 147107.859
$ws.Range("M126").Value = -144637.859

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 19616614
$ws.Range("I132").Value = 37049984
$ws.Range("J132").Value = 4074.75
$ws.Range("K132").Value = 111149952
$ws.Range("L132").Value = 12224.25
$ws.Range("M132").Value = -111147422
$ws.Range("N132").Value = -17284.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9814
$ws.Range("I136").Value = 7562.3335
$ws.Range("J136").Value = 14112.637
$ws.Range("K136").Value = 22687.0005
$ws.Range("L136").Value = 42337.911
$ws.Range("M136").Value = -20137.0005
$ws.Range("N136").Value = -47437.911

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 683.6129
$ws.Range("I126").Value = 420.77274
$ws.Range("K126").Value = 1262.31822
$ws.Range("M126").Value = 1207.68178

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1802.3143
$ws.Range("I132").Value = 1618.5
$ws.Range("J132").Value = 2203.3635
$ws.Range("K132").Value = 4855.5
$ws.Range("L132").Value = 6610.0905
$ws.Range("M132").Value = -2325.5
$ws.Range("N132").Value = -11670.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2519.2
$ws.Range("I136").Value = 3233.4119
$ws.Range("K136").Value = 9700.235700000001
$ws.Range("M136").Value = -7150.235700000001
